# BI-1613 - Updating test files
# Adds a new "Term Type" column (S) to the Template sheet and removes the
# stray bold/"applyFont" style that had been applied to the data rows
# (rows 3 and 4), including the now-empty N/O/P cells in those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column S: "Term Type" -------------------------------------------
$ws.Cells.Item(1, 19).Value = "Term Type"
$ws.Cells.Item(2, 19).Value = "Germplasm Attribute"
$ws.Cells.Item(4, 19).Value = "phenotype"
# (S3 intentionally left blank)

# --- Strip the extra style (s="1" / applyFont) from rows 3 & 4 -----------
$ws.Range("A3:D3").ClearFormats()
$ws.Range("F3:Q3").ClearFormats()
$ws.Range("A4:D4").ClearFormats()
$ws.Range("F4:Q4").ClearFormats()

# Those rows' N/O/P cells were empty placeholders that only existed to carry
# the style - remove them entirely now that the style is gone.
$ws.Range("N3:P3").ClearContents()
$ws.Range("N4:P4").ClearContents()

# --- View state: scroll so column Q is the leftmost visible column, and
# leave the selection on the new S5 cell (below the new column's data).
$ws.Range("S5").Select()
$excel.ActiveWindow.ScrollColumn = 17
